# Weekly status report push:
# 1) Update the "As Of" date
# 2) Update Percent Complete figure
# 3) Rewrite "Work Planned for Last Week" paragraph
# 4) Rewrite "Work Completed Last Week" paragraph
# 5) Rewrite "Work Planned for Next Week" paragraph
# 6) Touch "Open Issues" heading so any stale lastRenderedPageBreak marker clears

$d = $word.ActiveDocument

$wdReplaceAll = 2
$wdFindContinue = 1

# 1) Date: 03/28/2022 -> 04/03/2022
$d.Content.Find.Execute("03/28/2022", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "04/03/2022", $wdReplaceAll) | Out-Null

# 2) Percent Complete: 60% -> 65%
$d.Content.Find.Execute("60%", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "65%", $wdReplaceAll) | Out-Null

# 3) Work Planned for Last Week (old report -> this week's planned recap)
$oldPlannedLastWeek = "We planned to continue to refine the models. We will have to work with the sponsors to redefine the scope of the project or acquire a physical board from them that we can deploy the models to. As a result, both teams will continue to attempt achieving higher accuracies on the machine learning models while we wait to figure out how the scope will be redefined or see if we can get a physical board to use since emulation is not supported."
$newPlannedLastWeek = "Last week we planned to continue developing the quantum models and refine the classical models. The quantum SVM and quantum MLP are nearing completion. Additionally, the classical SVM and MLP are nearly complete as well. The only thing that remains is some fine tuning and then our analysis and comparison between the models. "
$d.Content.Find.Execute($oldPlannedLastWeek, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newPlannedLastWeek, $wdReplaceAll) | Out-Null

# 4) Work Completed Last Week
$oldCompleted = "We completed no work last week due to Spring Break."
$newCompleted = "We were able to get the quantum MLP working nicely. We still have some fine-tuning to do to get the accuracy above 90%. Additionally, the classical SVM and MLP are now completed and there are analysis notebooks for each. As we finish the quantum models we" + [char]8217 + "ll conduct more comparisons between the quantum and classical models. "
$d.Content.Find.Execute($oldCompleted, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newCompleted, $wdReplaceAll) | Out-Null

# 5) Work Planned for Next Week
$oldPlannedNextWeek = "We will continue to refine the models and work with the sponsors to redefine the scope of the project or acquire a physical board from them that we can deploy the models to. As a result, both teams will continue to attempt achieving higher accuracies on the machine learning models while we wait to figure out how the scope will be redefined or see if we can get a physical board to use since emulation is not supported. "
$newPlannedNextWeek = "This next week we are going to reach out to IBM to see if we can acquire some quantum tokens to deploy our models to their actual physical quantum computers. IBM has an API that will allow us to run our models on an actual quantum computer so we can get a realistic assessment of how well our quantum models are performing. We plan to reach out to IBM and figure out how to go about interfacing with their quantum API.  "
$d.Content.Find.Execute($oldPlannedNextWeek, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newPlannedNextWeek, $wdReplaceAll) | Out-Null

# 6) Normalize the "Open Issues" heading run (drops the stale lastRenderedPageBreak marker
#    left over from the previous pagination, same as Word does when the page reflows).
$d.Content.Find.Execute("Open Issues", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "Open Issues", $wdReplaceAll) | Out-Null

Write-Output "Status report updated."
